# Re-shuffle the taxon-identity columns (A, B, D, E, F, G, H, Q, R) across
# rows 8-35: for each destination row, pull the values that belong to a
# particular source row (per the authoritative mapping below), while every
# other column (location, dates, observers, etc.) stays put on its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow
$map = @{
    8  = 33
    9  = 35
    10 = 23
    11 = 15
    12 = 28
    13 = 8
    14 = 29
    15 = 34
    16 = 13
    17 = 12
    18 = 11
    19 = 22
    20 = 9
    21 = 18
    22 = 26
    23 = 17
    24 = 27
    25 = 30
    26 = 20
    27 = 10
    28 = 31
    29 = 32
    30 = 16
    31 = 14
    32 = 24
    33 = 21
    34 = 25
    35 = 19
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# First snapshot the source values for every row that participates, BEFORE
# any writes happen (the mapping is a full permutation, so writes would
# otherwise clobber values another destination still needs to read).
$snapshot = @{}
foreach ($row in $map.Keys) {
    $values = @{}
    foreach ($col in $cols) {
        $values[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $values
}

# Now write the snapshotted source row's values into each destination row.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $values = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $values[$col]
    }
}
